# Apply the "K instead of Strike#" regen: update column G (K) values for
# rows 2-13 on Sheet1 with the newly computed strikeout counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 0
    3  = 3
    4  = 2
    5  = 3
    6  = 4
    7  = 3
    8  = 2
    9  = 0
    10 = 1
    11 = 1
    12 = 2
    13 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
